$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 8689.714
$ws.Range("J40").Value = 8689.714
$ws.Range("L40").Value = 8689.714
$ws.Range("N40").Value = -9039.714

$ws.Range("H64").Value = 37044776
$ws.Range("I64").Value = 83337000
$ws.Range("J64").Value = 10998.6
$ws.Range("K64").Value = 83337000
$ws.Range("L64").Value = 10998.6
$ws.Range("M64").Value = -83336752
$ws.Range("N64").Value = -11494.6

$ws.Range("H67").Value = 37044776
$ws.Range("I67").Value = 83337000
$ws.Range("J67").Value = 10998.6
$ws.Range("K67").Value = 83337000
$ws.Range("L67").Value = 10998.6
$ws.Range("M67").Value = -83336142
$ws.Range("N67").Value = -12714.6

$ws.Range("H75").Value = 39999
$ws.Range("J75").Value = 39999
$ws.Range("L75").Value = 39999
$ws.Range("N75").Value = -41871

$ws.Range("H76").Value = 125004740
$ws.Range("I76").Value = 333336000
$ws.Range("K76").Value = 333336000
$ws.Range("M76").Value = -333335685

$ws.Range("H78").Value = 39999
$ws.Range("J78").Value = 39999
$ws.Range("L78").Value = 119997
$ws.Range("N78").Value = -129357

$ws.Range("H79").Value = 125004740
$ws.Range("I79").Value = 333336000
$ws.Range("K79").Value = 333336000
$ws.Range("M79").Value = -333334908

$ws.Range("H115").Value = 756.8570999999999
$ws.Range("I115").Value = 259.8
$ws.Range("J115").Value = 1999.5
$ws.Range("K115").Value = 779.4000000000001
$ws.Range("L115").Value = 5998.5
$ws.Range("M115").Value = 787.5999999999999
$ws.Range("N115").Value = -9132.5

$ws.Range("H118").Value = 591.5454999999999
$ws.Range("I118").Value = 191.33333
$ws.Range("J118").Value = 1071.8
$ws.Range("K118").Value = 573.99999
$ws.Range("L118").Value = 3215.4
$ws.Range("M118").Value = 1083.00001
$ws.Range("N118").Value = -6529.4

$ws.Range("H125").Value = 5414.2
$ws.Range("J125").Value = 9000
$ws.Range("L125").Value = 81000
$ws.Range("N125").Value = -85920

$ws.Range("H137").Value = 4360.1924
$ws.Range("I137").Value = 6174.3
$ws.Range("J137").Value = 3226.375
$ws.Range("K137").Value = 18522.9
$ws.Range("L137").Value = 9679.125
$ws.Range("M137").Value = -15972.9
$ws.Range("N137").Value = -14779.125

$ws.Range("H138").Value = 4512.472
$ws.Range("J138").Value = 4831.7856
$ws.Range("L138").Value = 14495.3568
$ws.Range("N138").Value = -24775.3568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 9999
$ws.Range("J46").Value = 9999
$ws.Range("L46").Value = 9999
$ws.Range("N46").Value = -10637

$ws.Range("H132").Value = 551446.5
$ws.Range("I132").Value = 585134.9
$ws.Range("J132").Value = 180875
$ws.Range("K132").Value = 1755404.7
$ws.Range("L132").Value = 542625
$ws.Range("M132").Value = -1752874.7
$ws.Range("N132").Value = -547685

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 7476.486
$ws.Range("I99").Value = 6955.326
$ws.Range("J99").Value = 8475.375
$ws.Range("K99").Value = 6955.326
$ws.Range("L99").Value = 8475.375
$ws.Range("M99").Value = -5457.326
$ws.Range("N99").Value = -11471.375

$ws.Range("H105").Value = 1942.7667
$ws.Range("I105").Value = 1839.3
$ws.Range("J105").Value = 2149.7
$ws.Range("K105").Value = 1839.3
$ws.Range("L105").Value = 2149.7
$ws.Range("M105").Value = -92.29999999999995
$ws.Range("N105").Value = -5643.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10877240
$ws.Range("I31").Value = 33347292
$ws.Range("J31").Value = 4634.484
$ws.Range("K31").Value = 33347292
$ws.Range("L31").Value = 4634.484
$ws.Range("M31").Value = -33346997
$ws.Range("N31").Value = -5224.484

$ws.Range("H34").Value = 10877240
$ws.Range("I34").Value = 33347292
$ws.Range("J34").Value = 4634.484
$ws.Range("K34").Value = 33347292
$ws.Range("L34").Value = 4634.484
$ws.Range("M34").Value = -33347090
$ws.Range("N34").Value = -5038.484

$ws.Range("H58").Value = 90924660
$ws.Range("I58").Value = 125009250
$ws.Range("K58").Value = 125009250
$ws.Range("M58").Value = -125009047

$ws.Range("H62").Value = 9655.714
$ws.Range("I62").Value = 10610.777
$ws.Range("J62").Value = 7936.6
$ws.Range("K62").Value = 10610.777
$ws.Range("L62").Value = 7936.6
$ws.Range("M62").Value = -9986.777
$ws.Range("N62").Value = -9184.6

$ws.Range("H65").Value = 9655.714
$ws.Range("I65").Value = 10610.777
$ws.Range("J65").Value = 7936.6
$ws.Range("K65").Value = 53053.885
$ws.Range("L65").Value = 39683
$ws.Range("M65").Value = -49933.885
$ws.Range("N65").Value = -45923

$ws.Range("H69").Value = 9500
$ws.Range("I69").Value = 9500
$ws.Range("K69").Value = 9500
$ws.Range("M69").Value = -8751

$ws.Range("H72").Value = 9500
$ws.Range("I72").Value = 9500
$ws.Range("K72").Value = 28500
$ws.Range("M72").Value = -24756

$ws.Range("H122").Value = 6733.875
$ws.Range("I122").Value = 2899.5557
$ws.Range("J122").Value = 11663.714
$ws.Range("K122").Value = 8698.667099999999
$ws.Range("L122").Value = 34991.142
$ws.Range("M122").Value = -6248.667099999999
$ws.Range("N122").Value = -39891.142

$ws.Range("H136").Value = 90924660
$ws.Range("I136").Value = 125009250
$ws.Range("K136").Value = 375027750
$ws.Range("M136").Value = -375025200

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 82.47059
$ws.Range("J23").Value = 96.72727
$ws.Range("L23").Value = 290.18181
$ws.Range("N23").Value = -760.18181

$ws.Range("H57").Value = 12797.6
$ws.Range("I57").Value = 7994
$ws.Range("K57").Value = 23982
$ws.Range("M57").Value = -23423

$ws.Range("H82").Value = 19165
$ws.Range("I82").Value = 7995.5
$ws.Range("J82").Value = 24749.75
$ws.Range("K82").Value = 23986.5
$ws.Range("L82").Value = 74249.25
$ws.Range("M82").Value = -23580.5
$ws.Range("N82").Value = -75061.25

$ws.Range("H85").Value = 19165
$ws.Range("I85").Value = 7995.5
$ws.Range("J85").Value = 24749.75
$ws.Range("K85").Value = 23986.5
$ws.Range("L85").Value = 74249.25
$ws.Range("M85").Value = -22582.5
$ws.Range("N85").Value = -77057.25

$ws.Range("H105").Value = 25933.334
$ws.Range("J105").Value = 25933.334
$ws.Range("L105").Value = 77800.00199999999
$ws.Range("N105").Value = -83042.00199999999

$ws.Range("H107").Value = 2632.7256
$ws.Range("J107").Value = 2999.182
$ws.Range("L107").Value = 8997.545999999998
$ws.Range("N107").Value = -12837.546

$ws.Range("H122").Value = 2961.0557
$ws.Range("I122").Value = 776.4
$ws.Range("J122").Value = 3801.3076
$ws.Range("K122").Value = 6987.599999999999
$ws.Range("L122").Value = 34211.7684
$ws.Range("M122").Value = -4537.599999999999
$ws.Range("N122").Value = -39111.7684

$ws.Range("H130").Value = 9999.75

$ws.Range("H134").Value = 18588
$ws.Range("I134").Value = 16605.6
$ws.Range("K134").Value = 49816.8
$ws.Range("M134").Value = -44746.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 30040.334
$ws.Range("I58").Value = 30040.334
$ws.Range("K58").Value = 30040.334
$ws.Range("M58").Value = -29763.334

$ws.Range("H70").Value = 5886.5625
$ws.Range("I70").Value = 5668.154
$ws.Range("J70").Value = 6833
$ws.Range("K70").Value = 5668.154
$ws.Range("L70").Value = 6833
$ws.Range("M70").Value = -5398.154
$ws.Range("N70").Value = -7373

$ws.Range("H73").Value = 5886.5625
$ws.Range("I73").Value = 5668.154
$ws.Range("J73").Value = 6833
$ws.Range("K73").Value = 5668.154
$ws.Range("L73").Value = 6833
$ws.Range("M73").Value = -4732.154
$ws.Range("N73").Value = -8705

$ws.Range("H80").Value = 9088.5
$ws.Range("I80").Value = 3302.5
$ws.Range("K80").Value = 3302.5
$ws.Range("M80").Value = -2304.5

$ws.Range("H83").Value = 9088.5
$ws.Range("I83").Value = 3302.5
$ws.Range("K83").Value = 16512.5
$ws.Range("M83").Value = -11520.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 911.7778
$ws.Range("J22").Value = 958.4
$ws.Range("L22").Value = 958.4
$ws.Range("N22").Value = -1548.4

$ws.Range("H27").Value = 911.7778
$ws.Range("J27").Value = 958.4
$ws.Range("L27").Value = 958.4
$ws.Range("N27").Value = -1172.4

$ws.Range("H40").Value = 5104.5
$ws.Range("I40").Value = 4735.615
$ws.Range("J40").Value = 9900
$ws.Range("K40").Value = 4735.615
$ws.Range("L40").Value = 9900
$ws.Range("M40").Value = -4599.615
$ws.Range("N40").Value = -10172

$ws.Range("H56").Value = 7054
$ws.Range("I56").Value = 4051
$ws.Range("K56").Value = 4051
$ws.Range("M56").Value = -3360

$ws.Range("H136").Value = 13801
$ws.Range("I136").Value = 16651.75
$ws.Range("K136").Value = 49955.25
$ws.Range("M136").Value = -47405.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 95750
$ws.Range("J46").Value = 95750
$ws.Range("L46").Value = 95750
$ws.Range("N46").Value = -96212

$ws.Range("H100").Value = 1949.7142
$ws.Range("I100").Value = 1678
$ws.Range("J100").Value = 2493.1428
$ws.Range("K100").Value = 3356
$ws.Range("L100").Value = 4986.2856
$ws.Range("M100").Value = -2815
$ws.Range("N100").Value = -6068.2856

$ws.Range("H103").Value = 32399.4
$ws.Range("J103").Value = 32399.4
$ws.Range("L103").Value = 32399.4
$ws.Range("N103").Value = -34743.4

$ws.Range("H107").Value = 8696424
$ws.Range("J107").Value = 972
$ws.Range("L107").Value = 2916
$ws.Range("N107").Value = -6756

$ws.Range("H113").Value = 6174100.5
$ws.Range("I113").Value = 10417903
$ws.Range("J113").Value = 1297.909
$ws.Range("K113").Value = 31253709
$ws.Range("L113").Value = 3893.727
$ws.Range("M113").Value = -31251539
$ws.Range("N113").Value = -8233.727000000001

$ws.Range("H134").Value = 95750
$ws.Range("J134").Value = 95750
$ws.Range("L134").Value = 287250
$ws.Range("N134").Value = -292320
